# Append: 2025-09-20 01:38 JST
# Update the "取得日時" (retrieved timestamp) column (A) for every data row
# on the first sheet ("ランサーズ") from the old run timestamp
# 2025-09-20 01:13:59 to the new one 2025-09-20 01:38:13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldTimestamp = "2025-09-20 01:13:59"
$newTimestamp = "2025-09-20 01:38:13"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
